$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 10 new rows (147-156) of MAC device data, continuing the pattern
# from the existing rows (regcntr_id, device_id, lang_code, is_active, cr_by, cr_dtimes)
$startRow = 147
$startDevice = 3000166

for ($i = 0; $i -lt 10; $i++) {
    $row = $startRow + $i
    $deviceId = $startDevice + $i

    $ws.Cells.Item($row, 1).Value = 10001
    $ws.Cells.Item($row, 2).Value = $deviceId
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
}

# Update the view to reflect scrolling to the new rows / selection state
$excel.ActiveWindow.ScrollRow = 140
$ws.Range("C152").Select()
